# Insert a new weekly price record as row 271 in the "Berenjena" sheet.
# This shifts the existing rows 271-283 down to 272-284 (dimension grows
# from A1:R283 to A1:R284) and populates the new row with the latest
# weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 271..283 down by one to make room for the new record.
$ws.Rows(271).Insert()

# Populate the newly inserted row 271 with the new weekly data.
$ws.Cells.Item(271, 1).Value = 9
$ws.Cells.Item(271, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(271, 3).Value = "Metropolitana"
$ws.Cells.Item(271, 4).Value = 44746
$ws.Cells.Item(271, 5).Value = 13
$ws.Cells.Item(271, 6).Value = 100112001
$ws.Cells.Item(271, 7).Value = "Berenjena"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 97
$ws.Cells.Item(271, 11).Value = 12000
$ws.Cells.Item(271, 12).Value = 13000
$ws.Cells.Item(271, 13).Value = 12495
$ws.Cells.Item(271, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(271, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(271, 16).Value = 250
$ws.Cells.Item(271, 17).Value = 50
$ws.Cells.Item(271, 18).Value = "Hortaliza"
